$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''52.005.76'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').Value = '''2.820.51'
$ws.Range('E3').Value = '  +2.99%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '''356.80'
$ws.Range('E5').Value = '  +7.42%  '
$ws.Range('D6').Value = '''113.93'
$ws.Range('E6').Value = '  -2.01%  '
$ws.Range('E7').Value = '  +2.83%  '
$ws.Range('D8').Value = '''0.999'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '''0.602'
$ws.Range('E9').Value = '  +5.25%  '
$ws.Range('D10').Value = '''41.80'
$ws.Range('E10').Value = '  +1.29%  '
$ws.Range('D11').Value = '''0.0847'
$ws.Range('E11').Value = '  +2.03%  '
$ws.Range('E12').Value = '  +1.49%  '
$ws.Range('D13').Value = '''19.88'
$ws.Range('E13').Value = '  -1.38%  '
$ws.Range('D14').Value = '''7.80'
$ws.Range('E14').Value = '  +3.37%  '
$ws.Range('D15').Value = '''3.249.70'
$ws.Range('E15').Value = '  +2.37%  '
$ws.Range('D16').Value = '''2.810.58'
$ws.Range('E16').Value = '  +2.35%  '
$ws.Range('E17').Value = '  +1.38%  '
$ws.Range('D18').Value = '''51.859.65'
$ws.Range('E18').Value = '  +0.84%  '
$ws.Range('B19').Value = 'ImmutableX'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D19').Value = '''3.17'
$ws.Range('E19').Value = '  +3.19%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '''7.27'
$ws.Range('E20').Value = '  +6.78%  '
$ws.Range('E21').Value = '  +1.45%  '
$ws.Range('D22').Value = '''0.0₃0984'
$ws.Range('E22').Value = '  +2.55%  '
$ws.Range('D23').Value = '''269.78'
$ws.Range('E23').Value = '  -2.68%  '
$ws.Range('D24').Value = '''69.55'
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('E25').Value = '  +6.43%  '
$ws.Range('D26').Value = '''26.80'
$ws.Range('E26').Value = '  +0.64%  '
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').Value = '''10.22'
$ws.Range('E28').Value = '  +0.54%  '
$ws.Range('E29').Value = '  +1.00%  '
$ws.Range('E30').Value = '  +0.58%  '
$ws.Range('D31').Value = '''50.64'
$ws.Range('E31').Value = '  +0.98%  '
$ws.Range('D32').Value = '''33.67'
$ws.Range('E32').Value = '  -3.50%  '
$ws.Range('D33').Value = '''5.84'
$ws.Range('E33').Value = '  +5.45%  '
$ws.Range('D34').Value = '''0.0436'
$ws.Range('E34').Value = '  +26.04%  '
$ws.Range('D35').Value = '''0.0827'
$ws.Range('E35').Value = '  +1.08%  '
$ws.Range('D36').Value = '''0.998'
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('D37').Value = '''2.09'
$ws.Range('E37').Value = '  +1.01%  '
$ws.Range('D38').Value = '''4.88'
$ws.Range('E38').Value = '  -0.46%  '
$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D39').Value = '''18.43'
$ws.Range('E39').Value = '  -3.05%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').Value = '''3.20'
$ws.Range('E40').Value = '  +1.67%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '''2.56'
$ws.Range('E41').Value = '  +6.19%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '''23.36'
$ws.Range('E42').Value = '  +1.27%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').Value = '''0.115'
$ws.Range('E43').Value = '  +1.80%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').Value = '''126.34'
$ws.Range('E44').Value = '  -2.32%  '
$ws.Range('D45').Value = '''2.28'
$ws.Range('E45').Value = '  +1.88%  '
$ws.Range('D46').Value = '''3.33'
$ws.Range('E46').Value = '  +0.27%  '
$ws.Range('D47').Value = '''2.045.06'
$ws.Range('E47').Value = '  -2.70%  '
$ws.Range('E48').Value = '  +3.87%  '
$ws.Range('B49').Value = 'SEI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range('D49').Value = '''0.945'
$ws.Range('E49').Value = '  +9.38%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = '''5.63'
$ws.Range('E50').Value = '  +2.51%  '
$ws.Range('D51').Value = '''8.90'
$ws.Range('E51').Value = '  -0.11%  '
